$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B29 was stored as an inline string "3"; change it to a real number 3
$ws.Cells.Item(29, 2).Value = 3

# Add new row 30 of annotation data
$ws.Cells.Item(30, 1).Value = "Ruilin"

# B30 must stay a text value "4" (not be auto-converted to a number)
$ws.Cells.Item(30, 2).NumberFormat = "@"
$ws.Cells.Item(30, 2).Value = "4"
$ws.Cells.Item(30, 2).Style = "Normal"

$ws.Cells.Item(30, 3).Value = "thank, detailed and insightful feedback"
$ws.Cells.Item(30, 4).Value = "ACK"
$ws.Cells.Item(30, 5).Value = "OTH"
$ws.Cells.Item(30, 6).Value = "7f314748-ac5a-4a11-8786-6125314f9d6d"
$ws.Cells.Item(30, 7).Value = "Sy2ogebAW_annotated.xlsx"
$ws.Cells.Item(30, 8).Value = "We would like to thank all reviewers for their detailed and insightful feedback."
